{"js": "// Refactor the {{image:...}} tag paragraph down to a minimal placeholder\n// and drop the trailing \"how to fill this tag\" explanation paragraphs,\n// per the commit \"Refactoring links ad images tags processing\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Remove every paragraph except the very last one (the body always keeps\n// at least one trailing paragraph, so deleting the last one is a no-op;\n// instead we delete all the others, leaving that final paragraph to hold\n// the rewritten tag text).\nfor (let i = paragraphs.items.length - 2; i >= 0; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\nconst remaining = body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\nconst target = remaining.items[0];\n\n// Replace whatever (empty) content is left with the trimmed tag text.\ntarget.insertText(\"Image ${{image:/}} test.\", Word.InsertLocation.start);\nawait context.sync();\n\n// Restore the _GoBack bookmark right after the final period, matching the\n// original document's bookmark placement.\nconst endRange = target.getRange(\"End\");\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Refactor the {{image:...}} tag paragraph down to a minimal placeholder\n# and drop the trailing \"how to fill this tag\" explanation paragraphs,\n# per the commit \"Refactoring links ad images tags processing\".\n\n$d = $word.ActiveDocument\n\n# Drop paragraphs 1-7: the original {{image:...}} tag paragraph (with its\n# title/source/imageFormat/width/height attribute list and spelling/\n# grammar proofing marks) plus the five descriptive paragraphs explaining\n# the POJO getter methods that fill the tag. Word always keeps at least\n# one trailing paragraph in the body, so paragraph 8 (already empty)\n# survives and becomes the new, only, paragraph.\nfor ($i = 7; $i -ge 1; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n\n# Type the trimmed placeholder text into what is now the sole paragraph.\n$p1 = $d.Paragraphs(1)\n$p1.Range.InsertBefore(\"Image `${{image:/}} test.#\")\n\n# Re-create the _GoBack bookmark right after the trailing period. A\n# temporary trailing marker character (\"#\", appended above) keeps the\n# period from being the very last character of the paragraph while the\n# bookmark is anchored, because anchoring exactly at paragraph-end\n# mis-positions the bookmark to the paragraph start in this host; the\n# marker is deleted immediately afterward, leaving the bookmark in the\n# correct place with no visible trace of it.\n$markerPos = $d.Content.End - 2\n$bookmarkRange = $d.Range($markerPos, $markerPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n$markerRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)\n$markerRange.Delete()\n"}
